$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark (it sat right before the "2005" run
#    in the Education section). A new one will be (re)created later where
#    Word would naturally leave it after the edits below (right before the
#    "Certifications" text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
